$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 2: J2/K2 ---
# Before: J2 = "сумма" (sum column header), K2 = "ТК"
# After:  J2 = "Лаб_1", K2 = empty
$ws.Range("J2").Value() = "Лаб_1"
$ws.Range("K2").ClearContents()

# --- Remove the "sum" formula column (J) and "TK" column (K) for all data rows ---
$ws.Range("J4:K32").ClearContents()

# --- Update scores for row 29 (Хромой Михаил) ---
$ws.Range("C29").Value() = 5
$ws.Range("D29").Value() = 5
$ws.Range("E29").Value() = 5
$ws.Range("F29").Value() = 5
$ws.Range("G29").Value() = 5
$ws.Range("H29").Value() = 5

# I29 is a brand-new cell (previously empty/nonexistent) that needs the same
# border/style formatting as the other "last homework column" cells (I16/I20).
$ws.Range("I16").Copy()
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("I29").Value() = 5
$excel.CutCopyMode = $false

# --- Frozen pane / scroll position & selection ---
$ws.Application.Goto($ws.Range("C10"), $false)
$ws.Range("J25").Select()
